$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")
$ws.Columns.Item(45).Delete()
